# Scheduled-runner style update of market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-crafting-class sheets. Values are plain data (no formulas), so each
# affected cell is written directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 713.6667
$ws.Range("I33").Value = 719
$ws.Range("K33").Value = 719
$ws.Range("M33").Value = -490

$ws.Range("H86").Value = 9810.556
$ws.Range("I86").Value = 3715.8333
$ws.Range("K86").Value = 3715.8333
$ws.Range("M86").Value = -2592.8333

$ws.Range("H89").Value = 9810.556
$ws.Range("I89").Value = 3715.8333
$ws.Range("K89").Value = 18579.1665
$ws.Range("M89").Value = -12963.1665

$ws.Range("H98").Value = 1606.3695
$ws.Range("I98").Value = 1639.5111
$ws.Range("K98").Value = 1639.5111
$ws.Range("M98").Value = -141.5110999999999

$ws.Range("H113").Value = 17545934
$ws.Range("I113").Value = 55557040
$ws.Range("J113").Value = 2345.8462
$ws.Range("K113").Value = 55557040
$ws.Range("L113").Value = 2345.8462
$ws.Range("M113").Value = -55553786
$ws.Range("N113").Value = -8853.8462

$ws.Range("H116").Value = 8318.091
$ws.Range("I116").Value = 6933.278
$ws.Range("J116").Value = 14549.75
$ws.Range("K116").Value = 6933.278
$ws.Range("L116").Value = 14549.75
$ws.Range("M116").Value = -3491.278
$ws.Range("N116").Value = -21433.75

$ws.Range("H122").Value = 1606.3695
$ws.Range("I122").Value = 1639.5111
$ws.Range("K122").Value = 4918.5333
$ws.Range("M122").Value = -2468.5333

$ws.Range("H132").Value = 6091.5405
$ws.Range("I132").Value = 6670.4194
$ws.Range("J132").Value = 3100.6667
$ws.Range("K132").Value = 20011.2582
$ws.Range("L132").Value = 9302.000100000001
$ws.Range("M132").Value = -17481.2582
$ws.Range("N132").Value = -14362.0001

$ws.Range("H135").Value = 774.675
$ws.Range("I135").Value = 649.34283
$ws.Range("K135").Value = 5844.08547
$ws.Range("M135").Value = -3309.08547

$ws.Range("H137").Value = 1518480.5
$ws.Range("I137").Value = 6252751.5
$ws.Range("J137").Value = 3513.72
$ws.Range("K137").Value = 18758254.5
$ws.Range("L137").Value = 10541.16
$ws.Range("M137").Value = -18755704.5
$ws.Range("N137").Value = -15641.16

$ws.Range("H138").Value = 3496.132
$ws.Range("I138").Value = 2996.3044
$ws.Range("K138").Value = 8988.913199999999
$ws.Range("M138").Value = -3848.913199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2477.6206
$ws.Range("I61").Value = 2221.9524
$ws.Range("K61").Value = 2221.9524
$ws.Range("M61").Value = -2009.9524

$ws.Range("H74").Value = 329714.12
$ws.Range("I74").Value = 465345
$ws.Range("K74").Value = 465345
$ws.Range("M74").Value = -464471

$ws.Range("H77").Value = 329714.12
$ws.Range("I77").Value = 465345
$ws.Range("K77").Value = 2326725
$ws.Range("M77").Value = -2322357

$ws.Range("H122").Value = 5643.4346
$ws.Range("I122").Value = 5333.278
$ws.Range("J122").Value = 6760
$ws.Range("K122").Value = 15999.834
$ws.Range("L122").Value = 20280
$ws.Range("M122").Value = -13549.834
$ws.Range("N122").Value = -25180

$ws.Range("H136").Value = 2477.6206
$ws.Range("I136").Value = 2221.9524
$ws.Range("K136").Value = 6665.8572
$ws.Range("M136").Value = -4115.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16135236
$ws.Range("I20").Value = 22734460
$ws.Range("K20").Value = 22734460
$ws.Range("M20").Value = -22734213

$ws.Range("H86").Value = 3099.889
$ws.Range("I86").Value = 2628.4285
$ws.Range("K86").Value = 2628.4285
$ws.Range("M86").Value = -1505.4285

$ws.Range("H89").Value = 3099.889
$ws.Range("I89").Value = 2628.4285
$ws.Range("K89").Value = 13142.1425
$ws.Range("M89").Value = -7526.1425

$ws.Range("H99").Value = 46426.26
$ws.Range("I99").Value = 57669.332
$ws.Range("K99").Value = 57669.332
$ws.Range("M99").Value = -56171.332

$ws.Range("H134").Value = 3226.3242
$ws.Range("I134").Value = 2774.8462
$ws.Range("K134").Value = 8324.5386
$ws.Range("M134").Value = -5789.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4039590
$ws.Range("I31").Value = 6582.6665
$ws.Range("J31").Value = 9623754
$ws.Range("K31").Value = 6582.6665
$ws.Range("L31").Value = 9623754
$ws.Range("M31").Value = -6287.6665
$ws.Range("N31").Value = -9624344

$ws.Range("H34").Value = 4039590
$ws.Range("I34").Value = 6582.6665
$ws.Range("J34").Value = 9623754
$ws.Range("K34").Value = 6582.6665
$ws.Range("L34").Value = 9623754
$ws.Range("M34").Value = -6380.6665
$ws.Range("N34").Value = -9624158

$ws.Range("H107").Value = 3572680.8
$ws.Range("J107").Value = 1601.8572
$ws.Range("L107").Value = 1601.8572
$ws.Range("N107").Value = -5441.8572

$ws.Range("H120").Value = 35810.5
$ws.Range("I120").Value = 34256
$ws.Range("K120").Value = 34256
$ws.Range("M120").Value = -30627

$ws.Range("H132").Value = 13022512
$ws.Range("I132").Value = 9616543
$ws.Range("K132").Value = 28849629
$ws.Range("M132").Value = -28847099

$ws.Range("H141").Value = 312499.5
$ws.Range("J141").Value = 312499.5
$ws.Range("L141").Value = 312499.5
$ws.Range("N141").Value = -322859.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1874.3077
$ws.Range("I2").Value = 264.25
$ws.Range("J2").Value = 2589.889
$ws.Range("K2").Value = 1585.5
$ws.Range("L2").Value = 15539.334
$ws.Range("M2").Value = -1472.5
$ws.Range("N2").Value = -15765.334

$ws.Range("H5").Value = 1024.4286
$ws.Range("I5").Value = 1044.25
$ws.Range("J5").Value = 998
$ws.Range("K5").Value = 3132.75
$ws.Range("L5").Value = 2994
$ws.Range("M5").Value = -3020.75
$ws.Range("N5").Value = -3218

$ws.Range("H75").Value = 1255
$ws.Range("J75").Value = 1483.4
$ws.Range("L75").Value = 4450.200000000001
$ws.Range("N75").Value = -6446.200000000001

$ws.Range("H78").Value = 1255
$ws.Range("J78").Value = 1483.4
$ws.Range("L78").Value = 13350.6
$ws.Range("N78").Value = -23334.6

$ws.Range("H98").Value = 1785
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1785
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 5355
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -8351

$ws.Range("H113").Value = 1156.2307
$ws.Range("I113").Value = 273.33334
$ws.Range("J113").Value = 1271.3914
$ws.Range("K113").Value = 820.0000200000001
$ws.Range("L113").Value = 3814.1742
$ws.Range("M113").Value = 1349.99998
$ws.Range("N113").Value = -8154.174199999999

$ws.Range("H132").Value = 3000
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060

$ws.Range("H135").Value = 1024.4286
$ws.Range("I135").Value = 1044.25
$ws.Range("J135").Value = 998
$ws.Range("K135").Value = 9398.25
$ws.Range("L135").Value = 8982
$ws.Range("M135").Value = -6863.25
$ws.Range("N135").Value = -14052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 10133.577
$ws.Range("J24").Value = 10061.667
$ws.Range("L24").Value = 10061.667
$ws.Range("N24").Value = -10407.667

$ws.Range("H46").Value = 29918.2
$ws.Range("J46").Value = 29918.2
$ws.Range("L46").Value = 29918.2
$ws.Range("N46").Value = -30230.2

$ws.Range("H70").Value = 136159.31
$ws.Range("I70").Value = 253331.25
$ws.Range("J70").Value = 18987.375
$ws.Range("K70").Value = 253331.25
$ws.Range("L70").Value = 18987.375
$ws.Range("M70").Value = -253061.25
$ws.Range("N70").Value = -19527.375

$ws.Range("H73").Value = 136159.31
$ws.Range("I73").Value = 253331.25
$ws.Range("J73").Value = 18987.375
$ws.Range("K73").Value = 253331.25
$ws.Range("L73").Value = 18987.375
$ws.Range("M73").Value = -252395.25
$ws.Range("N73").Value = -20859.375

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 2541.8438
$ws.Range("I132").Value = 1613.8572
$ws.Range("J132").Value = 4313.4546
$ws.Range("K132").Value = 4841.571599999999
$ws.Range("L132").Value = 12940.3638
$ws.Range("M132").Value = -2311.571599999999
$ws.Range("N132").Value = -18000.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 263.66666
$ws.Range("I55").Value = 194
$ws.Range("K55").Value = 194
$ws.Range("M55").Value = -21

$ws.Range("H102").Value = 40999.332
$ws.Range("J102").Value = 40999.332
$ws.Range("L102").Value = 40999.332
$ws.Range("N102").Value = -47489.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19983
$ws.Range("J41").Value = 19983
$ws.Range("L41").Value = 19983
$ws.Range("N41").Value = -20763

$ws.Range("H96").Value = 1430
$ws.Range("I96").Value = 1401.6666
$ws.Range("K96").Value = 1401.6666
$ws.Range("M96").Value = -28.66660000000002

$ws.Range("H122").Value = 7144411
$ws.Range("I122").Value = 1256.7333
$ws.Range("K122").Value = 3770.199900000001
$ws.Range("M122").Value = -1320.199900000001

$ws.Range("H132").Value = 12824769
$ws.Range("I132").Value = 22226780
$ws.Range("K132").Value = 66680340
$ws.Range("M132").Value = -66677810

$ws.Range("H136").Value = 5938.1953
$ws.Range("I136").Value = 5846.355
$ws.Range("K136").Value = 17539.065
$ws.Range("M136").Value = -14989.065
